$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Status text: "Ready for handoff" -> "Handed back: in sync with en-US"
#    (Overview sheet status columns per-language, and each language sheet's
#    own Status column)
# ---------------------------------------------------------------------------
$statusText = "Handed back: in sync with en-US"

$ov = $wb.Worksheets.Item("Overview")
$ov.Range("E2").Value = $statusText
$ov.Range("F2").Value = $statusText
$ov.Range("E3").Value = $statusText
$ov.Range("F3").Value = $statusText

$zh = $wb.Worksheets.Item("zh-cn")
$zh.Range("C2").Value = $statusText
$zh.Range("C3").Value = $statusText

$de = $wb.Worksheets.Item("de-de")
$de.Range("C2").Value = $statusText
$de.Range("C3").Value = $statusText

# ---------------------------------------------------------------------------
# 2. Populate "Latest Target File" (I) / "Latest Handback File" (J) /
#    "Latest Handback DateTime" (K) for both data rows on the zh-cn and
#    de-de sheets - this is the handback report data written once the
#    handback package has been generated.
# ---------------------------------------------------------------------------

$baseUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/dd1223fd89fa59740dbb538600c1c1c955fd2e16/e2e/"

# --- zh-cn sheet ---
$zh.Range("J2").Value = "58dc30ff-4c75-47e5-95bf-406b1a6b723e.f5be56a1146a15bd60a4a7cfdf825f8351aff8e5.zh-cn.xlf"
$zh.Range("K2").Value = "2016-08-21 03:04:41"
$zh.Hyperlinks.Add($zh.Range("I2"), ($baseUrl + "58dc30ff-4c75-47e5-95bf-406b1a6b723e.md"), [Type]::Missing, [Type]::Missing, "58dc30ff-4c75-47e5-95bf-406b1a6b723e.md") | Out-Null

$zh.Range("J3").Value = "5ec2340e-2472-47a6-8ec5-01d3c83f9145.8a631ac5f48729918fecb1585a2f74ac01b3fb22.zh-cn.xlf"
$zh.Range("K3").Value = "2016-08-21 03:04:41"
$zh.Hyperlinks.Add($zh.Range("I3"), ($baseUrl + "5ec2340e-2472-47a6-8ec5-01d3c83f9145.md"), [Type]::Missing, [Type]::Missing, "5ec2340e-2472-47a6-8ec5-01d3c83f9145.md") | Out-Null

# --- de-de sheet ---
$de.Range("J2").Value = "58dc30ff-4c75-47e5-95bf-406b1a6b723e.f5be56a1146a15bd60a4a7cfdf825f8351aff8e5.de-de.xlf"
$de.Range("K2").Value = "2016-08-21 03:04:48"
$de.Hyperlinks.Add($de.Range("I2"), ($baseUrl + "58dc30ff-4c75-47e5-95bf-406b1a6b723e.md"), [Type]::Missing, [Type]::Missing, "58dc30ff-4c75-47e5-95bf-406b1a6b723e.md") | Out-Null

$de.Range("J3").Value = "5ec2340e-2472-47a6-8ec5-01d3c83f9145.8a631ac5f48729918fecb1585a2f74ac01b3fb22.de-de.xlf"
$de.Range("K3").Value = "2016-08-21 03:04:48"
$de.Hyperlinks.Add($de.Range("I3"), ($baseUrl + "5ec2340e-2472-47a6-8ec5-01d3c83f9145.md"), [Type]::Missing, [Type]::Missing, "5ec2340e-2472-47a6-8ec5-01d3c83f9145.md") | Out-Null

# ---------------------------------------------------------------------------
# 3. Widen the columns whose content just grew (Status columns on all
#    sheets, and the new Target/Handback File columns) so the report is
#    still readable, mirroring the wider columns the report generator
#    produces once a status/filename column holds longer text.
# ---------------------------------------------------------------------------
$ov.Range("E1").ColumnWidth = 29.083333333333332
$ov.Range("F1").ColumnWidth = 29.083333333333332

$zh.Range("C1").ColumnWidth = 29.083333333333332
$zh.Range("I1").ColumnWidth = 39.083333333333336
$zh.Range("J1").ColumnWidth = 39.083333333333336

$de.Range("C1").ColumnWidth = 29.083333333333332
$de.Range("I1").ColumnWidth = 39.083333333333336
$de.Range("J1").ColumnWidth = 39.083333333333336
